# Add a new "service" entry to the CV entries worksheet:
#   Core Member (two-year appointment): Social Systems Data Science Network
# This entry belongs right before the existing "Faculty Advisory Committee..."
# service row, so we insert a new row at row 15 and push everything else down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 15 (shifts rows 15-34 down to 16-35).
# Excel copies the row-above's cell formatting into the newly inserted row,
# which is why B15/C15/D15 below already carry the s="2"/s="3"/s="1" styles
# before we touch them.
$ws.Rows.Item(15).Insert()

# Populate the new row 15 with the new service entry.
$ws.Range("A15").Value = "service"
$ws.Range("B15").Value = "2019"
$ws.Range("C15").Value = "current"
$ws.Range("D15").WrapText = $true
$ws.Range("D15").Value = "Core Member (two-year appointment): Social Systems Data Science Network"
$ws.Range("E15").Value = "University of Oregon"

# The two-line wrapped text needs the taller row height the other wrapped
# "service" rows use (matches the saved row height for this entry).
$ws.Rows.Item(15).RowHeight = 34

# Match the author's final selection/scroll position from the saved workbook.
$ws.Range("E16").Select()
